# Add a 10th value to the color palette on the "simulation" sheet.
# - rows 2-10 get a new (matplotlib "tab10") color value
# - a new row 11 is added for "T cell CD4+" with the 10th color
# - the B column cells touched get a taller row (14.9) and wrap text enabled
# - the last selected cell becomes F14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("simulation")
$ws.Activate()

# New 10-color palette (matplotlib tab10) replacing the old 8-color one.
$colors = @{
    2  = "#1f77b4"
    3  = "#ff7f0e"
    4  = "#2ca02c"
    5  = "#d62728"
    6  = "#9467bd"
    7  = "#8c564b"
    8  = "#e377c2"
    9  = "#7f7f7f"
    10 = "#bcbd22"
    11 = "#17becf"
}

# Row 11 is brand new and needs its label in column A.
$ws.Range("A11").Value = "T cell CD4+"

foreach ($row in 2..11) {
    $ws.Range("B$row").Value = $colors[$row]
}

# Rows 2 and 4-11 grow a bit taller and the color cells wrap their text;
# row 3 is left untouched (its color only changes because the shared
# string it points at changes meaning).
foreach ($row in @(2,4,5,6,7,8,9,10,11)) {
    $ws.Rows.Item($row).RowHeight = 14.9
    $ws.Range("B$row").WrapText = $true
}

# Reflect the final cursor position recorded in the saved file.
$ws.Range("F14").Select()
